$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Офлайн режим работы"
$ws.Range("A12").Select()
